# Commit: "Math and test libraries reworked"
#
# The "Больничные" (Sick-leave) sheet had a bug in the sick-pay-per-day
# array formula in column H: the second INDEX() into Сотрудники!K:K used
# Больничные!B<row> instead of Больничные!B<row>+1 (an off-by-one versus
# the matching lookup already used in the first INDEX()). Fixing the
# lookup ripples through I/J on that sheet and P/R on "Сотрудники".
#
# This script only touches the actual content (formulas) and the saved
# cursor/selection per sheet, matching the authoritative diff.

$wb  = $excel.ActiveWorkbook

$wsEmployees = $wb.Worksheets.Item("Сотрудники")
$wsFines     = $wb.Worksheets.Item("Штрафы")
$wsSick      = $wb.Worksheets.Item("Больничные")

# --- "Больничные" sheet: fix the off-by-one in the H column array formula ---
# Old: INDEX(Сотрудники!H:H,Больничные!B<r>+1)*INDEX(Сотрудники!K:K,Больничные!B<r>)
# New: INDEX(Сотрудники!H:H,Больничные!B<r>+1)*INDEX(Сотрудники!K:K,Больничные!B<r>+1)
$wsSick.Range("H2").FormulaArray = "=INDEX(Сотрудники!H:H,Больничные!B2+1)*INDEX(Сотрудники!K:K,Больничные!B2+1)"
$wsSick.Range("H3").FormulaArray = "=INDEX(Сотрудники!H:H,Больничные!B3+1)*INDEX(Сотрудники!K:K,Больничные!B3+1)"
$wsSick.Range("H4").FormulaArray = "=INDEX(Сотрудники!H:H,Больничные!B4+1)*INDEX(Сотрудники!K:K,Больничные!B4+1)"

# P4/R4 on "Сотрудники" are driven by SUMIFS()/shared formulas that read
# the corrected J column on "Больничные" above - they recompute on their
# own, no direct edit needed there.

# --- Restore each sheet's saved selection (per-sheet cursor position) ---
# Selecting on a non-active sheet activates it, so set the other sheets
# first and re-activate + reselect "Сотрудники" (the tab that should stay
# selected) last.
$wsFines.Range("F2").Select() | Out-Null
$wsSick.Range("I10").Select() | Out-Null

$wsEmployees.Activate() | Out-Null
$wsEmployees.Range("Q14").Select() | Out-Null
